$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old extent (previous table used columns A:G) before writing new data.
# ClearContents() removes values only, so the pre-existing cell formatting (the
# date-style format carried on column C) is preserved for the rows that already
# existed (rows 2 and 3).
$ws.Range("A1:G6").ClearContents()

# Header row
$ws.Range("A1").Value = "Test Case ID"
$ws.Range("B1").Value = "Execution"
$ws.Range("C1").Value = "Class"
$ws.Range("D1").Value = "Browser"
$ws.Range("E1").Value = "Suite"

# Data rows
$ws.Range("A2").Value = "TC_004"
$ws.Range("B2").Value = "N"
$ws.Range("C2").Value = "LoginTest"
$ws.Range("D2").Value = "chrome"
$ws.Range("E2").Value = "Regression"

$ws.Range("A3").Value = "TC_005"
$ws.Range("B3").Value = "N"
$ws.Range("C3").Value = "LoginTest"
$ws.Range("D3").Value = "chrome"
$ws.Range("E3").Value = "Smoke"

$ws.Range("A4").Value = "TC_004"
$ws.Range("B4").Value = "N"
$ws.Range("C4").Value = "LoginTest"
$ws.Range("D4").Value = "firefox"
$ws.Range("E4").Value = "Smoke"

$ws.Range("A5").Value = "TC_013"
$ws.Range("B5").Value = "Y"
$ws.Range("C5").Value = "CartTest"
$ws.Range("D5").Value = "chrome"
$ws.Range("E5").Value = "Regression"

$ws.Range("A6").Value = "TC_015"
$ws.Range("B6").Value = "Y"
$ws.Range("C6").Value = "CartTest"
$ws.Range("D6").Value = "firefox"
$ws.Range("E6").Value = "Regression"

# Propagate the original "Class" column date-style formatting onto the new rows
# (4-6) that did not exist before, so the whole column is styled consistently.
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C6").PasteSpecial(-4122)

$ws.Range("B5").Select()
